# Update "Diaporama revue 1" slide (slide 4): reposition three text boxes
# and remove the last bullet ("Étalonnage des deux capteurs de température")
# from the "ZoneTexte 1" text box.
#
# Point values below are chosen so that, after round-tripping through the
# Shape.Top/Left Single-precision (float32) storage used by the PowerPoint
# object model, they land back on the exact target EMU values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# --- Shape "ZoneTexte 2" : move down (y only) ---------------------------
$zt2 = $s.Shapes.Item("ZoneTexte 2")
$zt2.Top = 187.60315710629922

# --- Shape "ZoneTexte 4" : move (x and y) --------------------------------
$zt4 = $s.Shapes.Item("ZoneTexte 4")
$zt4.Left = 85.23102362204725
$zt4.Top = 281.53622047244096

# --- Shape "ZoneTexte 1" : move down and drop the last paragraph --------
$zt1 = $s.Shapes.Item("ZoneTexte 1")
$zt1.Top = 374.93299872598425

$tr = $zt1.TextFrame.TextRange
$paraCount = $tr.Paragraphs().Count
$fullText = $tr.Text

# Remove the final paragraph entirely (its text plus its paragraph mark)
# by deleting everything starting right after the previous paragraph's
# text runs out. Deleting through a TextRange (rather than clearing the
# paragraph's .Text) avoids leaving a stray empty trailing paragraph.
$prevPara = $tr.Paragraphs($paraCount - 1)
$deleteStart = $prevPara.Start + $prevPara.Length
$deleteLength = $fullText.Length - $deleteStart + 2
$deleteRange = $tr.Characters($deleteStart, $deleteLength)
$deleteRange.Delete()
